$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy column D (rows 1-208) formats+values into column E, mirroring the
#    new "HeroConfigIDEx" column added alongside "DefaultObject".
$ws.Range("D1:D208").Copy()
$ws.Range("E1:E208").PasteSpecial(-4104)  # xlPasteAll

# 2. Header text for the new column.
$ws.Range("E1").Value = "HeroConfigIDEx"

# 3. For rows 158-208 the new column holds "Hero001" instead of "DefaultObject".
$ws.Range("E158:E208").Value = "Hero001"

# 4. Column E should be the same width as column D.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# 5. Update the list-validation range: it used to skip column D rows 7:9 via E:H,
#    now D8:E9 are included together and F:H keep the list validation.
$ws.Range("B7:C9,D8:E9,F7:H9").Validation.Delete()
$ws.Range("B7:C9,D8:E9,F7:H9").Validation.Add(3, 1, 1, "TRUE,FALSE")

# 6. Update view state: scrolled/selected position after editing.
$ws.Application.ActiveWindow.ScrollRow = 166
$ws.Range("E208").Select()
